{"js": "const replacements = [\n  [\"112\u00f76=18, 4\", \"284\u00f77=40, 4\"],\n  [\"741\u00f76=123, 3\", \"177\u00f77=25, 2\"],\n  [\"550\u00f73=183, 1\", \"245\u00f77=35, 0\"],\n  [\"558\u00f74=139, 2\", \"577\u00f72=288, 1\"],\n  [\"372\u00f79=41, 3\", \"389\u00f73=129, 2\"],\n  [\"296\u00f72=148, 0\", \"699\u00f78=87, 3\"],\n  [\"836\u00f79=92, 8\", \"894\u00f73=298, 0\"],\n  [\"657\u00f78=82, 1\", \"444\u00f79=49, 3\"],\n  [\"398\u00f75=79, 3\", \"666\u00f76=111, 0\"],\n  [\"704\u00f78=88, 0\", \"182\u00f79=20, 2\"],\n  [\"994\u00f78=124, 2\", \"432\u00f76=72, 0\"],\n  [\"532\u00f73=177, 1\", \"691\u00f76=115, 1\"],\n  [\"638\u00f79=70, 8\", \"108\u00f76=18, 0\"],\n  [\"452\u00f74=113, 0\", \"595\u00f72=297, 1\"],\n  [\"756\u00f78=94, 4\", \"878\u00f79=97, 5\"],\n  [\"465\u00f73=155, 0\", \"291\u00f74=72, 3\"],\n  [\"710\u00f76=118, 2\", \"385\u00f75=77, 0\"],\n  [\"753\u00f73=251, 0\", \"548\u00f75=109, 3\"],\n  [\"962\u00f77=137, 3\", \"378\u00f78=47, 2\"],\n  [\"954\u00f78=119, 2\", \"172\u00f79=19, 1\"],\n  [\"778\u00f75=155, 3\", \"770\u00f77=110, 0\"],\n  [\"213\u00f76=35, 3\", \"805\u00f73=268, 1\"],\n  [\"961\u00f73=320, 1\", \"165\u00f72=82, 1\"],\n  [\"836\u00f76=139, 2\", \"743\u00f77=106, 1\"],\n  [\"416\u00f75=83, 1\", \"245\u00f72=122, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"112\u00f76=18, 4\", \"284\u00f77=40, 4\"),\n    @(\"741\u00f76=123, 3\", \"177\u00f77=25, 2\"),\n    @(\"550\u00f73=183, 1\", \"245\u00f77=35, 0\"),\n    @(\"558\u00f74=139, 2\", \"577\u00f72=288, 1\"),\n    @(\"372\u00f79=41, 3\", \"389\u00f73=129, 2\"),\n    @(\"296\u00f72=148, 0\", \"699\u00f78=87, 3\"),\n    @(\"836\u00f79=92, 8\", \"894\u00f73=298, 0\"),\n    @(\"657\u00f78=82, 1\", \"444\u00f79=49, 3\"),\n    @(\"398\u00f75=79, 3\", \"666\u00f76=111, 0\"),\n    @(\"704\u00f78=88, 0\", \"182\u00f79=20, 2\"),\n    @(\"994\u00f78=124, 2\", \"432\u00f76=72, 0\"),\n    @(\"532\u00f73=177, 1\", \"691\u00f76=115, 1\"),\n    @(\"638\u00f79=70, 8\", \"108\u00f76=18, 0\"),\n    @(\"452\u00f74=113, 0\", \"595\u00f72=297, 1\"),\n    @(\"756\u00f78=94, 4\", \"878\u00f79=97, 5\"),\n    @(\"465\u00f73=155, 0\", \"291\u00f74=72, 3\"),\n    @(\"710\u00f76=118, 2\", \"385\u00f75=77, 0\"),\n    @(\"753\u00f73=251, 0\", \"548\u00f75=109, 3\"),\n    @(\"962\u00f77=137, 3\", \"378\u00f78=47, 2\"),\n    @(\"954\u00f78=119, 2\", \"172\u00f79=19, 1\"),\n    @(\"778\u00f75=155, 3\", \"770\u00f77=110, 0\"),\n    @(\"213\u00f76=35, 3\", \"805\u00f73=268, 1\"),\n    @(\"961\u00f73=320, 1\", \"165\u00f72=82, 1\"),\n    @(\"836\u00f76=139, 2\", \"743\u00f77=106, 1\"),\n    @(\"416\u00f75=83, 1\", \"245\u00f72=122, 1\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}"}
